$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Read (Academic)")

# A2: date-like text must stay literal text "03/12/24" (not an Excel date serial).
# Leading apostrophe forces text entry; resetting the style back to Normal afterwards
# clears the "quote prefix" formatting Excel applies to ambiguous text so the cell
# ends up with no explicit style, same as the other plain text cells.
$ws.Cells.Item(2, 1).Value = "'03/12/24"
$ws.Cells.Item(2, 1).Style = "Normal"

# C2 / D2: plain sign-in / sign-out times as literal text.
$ws.Cells.Item(2, 3).Value = "12:41"
$ws.Cells.Item(2, 4).Value = "12:48"

# E2: elapsed time (12:48 - 12:41 = 7 minutes) stored as a numeric day-fraction,
# formatted with the sheet's existing [hh]:mm:ss duration format.
$ws.Cells.Item(2, 5).Style = "Normal"
$ws.Cells.Item(2, 5).NumberFormat = "[hh]:mm:ss"
$ws.Cells.Item(2, 5).Value = 0.004861111111111111

# Grow the table to include the newly populated row (A1:E2 -> A1:E3) so a
# fresh blank "insert row" is available below the data, matching how Excel
# extends a table when its insert row is filled in.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E3"))
